# Apply the CSV-import template update:
#  - shared string "name of department 1/name of department 2"
#    becomes "name of department 1|name of department 2"
#  - the active sheet's selection moves from F4 to E8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the cells that hold the department-name placeholder text so the
# shared string itself changes from a "/" separator to a "|" separator.
$ws.Range("F2:F4").Value = "name of department 1|name of department 2"

# Move the selection/active cell to E8.
$ws.Range("E8").Select()
